# Rename the "MS" worksheet to "MS_Heat" so the sheet tab name matches
# the assay's folder name (assays/MS_Heat/isa.assay.xlsx).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MS")
$ws.Name = "MS_Heat"
